$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (the "Förändrad" date column), rows 2 through 260, all hold the
# same date serial value (45205 = 2023-10-06) that needs to be advanced by
# one day to 45206 (2023-10-07) for every data row.
$lastRow = 260
$range = $ws.Range("C2:C$lastRow")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
